$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the old "Lab 1 Descriptives" lab into a new video-based "Basics" lab,
# and rename the old "Lab 2 Variance and Transformed scores" lab to "Lab 2 Descriptives".
$ws.Range("C3").Value = "Lab 1 Basics"
$ws.Range("C4").Value = "Lab 2 Descriptives"

# Leave the selection where the editor ended up after making the change.
$ws.Range("C5").Select()
